# Final_Report.xlsx / "Summary" sheet:
#  - widen the report from a 3-column (Account | Amount doc curr | Amount local
#    curr) layout to a 6-column layout that also shows Company, Document
#    currency and Local Currency, and bump the report date.
#  - debugged duplication and fixed api issue

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------------
# Row 1 - title bar: extend the merged banner from A1:C1 to A1:F1. Copy the
# existing blank-but-styled C1 cell's formatting onto the three new trailing
# cells before re-merging (A1 already holds the title and keeps it because it
# is the top-left cell of the merge).
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A1:F1").Merge()

# ---------------------------------------------------------------------------
# Row 2 - report date stays in B2 (not shifted); only the text changes.
# Force literal text so Excel doesn't coerce it into a real date serial.
# ---------------------------------------------------------------------------
$ws.Range("B2").Formula = "'2025-09-04"

# ---------------------------------------------------------------------------
# Row 4 - header row. Old layout was Account | Amount in doc. curr. | Amount
# in local currency (A/B/C). New layout inserts Company before Account,
# Document currency before Amount in doc. curr., and Local Currency before
# Amount in local currency, i.e. old A->B, old B->D, old C->F.
# Copy the existing header style (A4, s=1: bold/centered/bordered) across the
# whole row first, then overwrite the text.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A4:F4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A4").Value2 = "Comapany"
$ws.Range("B4").Value2 = "Account"
$ws.Range("C4").Value2 = "Document currency"
$ws.Range("D4").Value2 = "Amount in doc. curr."
$ws.Range("E4").Value2 = "Local Currency"
$ws.Range("F4").Value2 = "Amount in local currency"

# ---------------------------------------------------------------------------
# Rows 5-24 - data. Same column remap as the header: old A->B (account),
# old B->D (amount in doc. curr.), old C->F (amount in local currency), with
# new Company/Document currency/Local Currency columns filled in (A/C/E).
# Account numbers are numeric-looking strings, so they are entered with a
# leading apostrophe to keep them as text instead of being parsed as numbers.
# ---------------------------------------------------------------------------
$data = @(
    @("UN0100", "63010001", "USD", -6850603.8, "USD", -6850603.8),
    @("XT0150", "63010002", "USD", 11501418.04, "USD", 11501418.04),
    @("XT0150", "63010011", "LKR", -1479775595.78, "USD", -7946810.529999999),
    @("UN0150", "63010012", "LKR", 1495088917.7, "USD", 8028294.79),
    @("XT0150", "63010061", "HKD", -1800000, "USD", -232183.17),
    @("XT0150", "63010101", "EUR", -41779.00999999999, "USD", -46504.21),
    @("UN0100", "63010162", "USD", 1354247.56, "USD", 1354247.56),
    @("UN0150", "63010502", "USD", 1368805.95, "USD", 1368805.95),
    @("XT0150", "63011001", "USD", -506397.46, "USD", -506397.46),
    @("XT0150", "63011011", "LKR", -181229989.54, "USD", -973255.86),
    @("XT0151", "63020001", "LKR", -3354998.88, "USD", -18017.3),
    @("XT0150", "63020002", "LKR", 3533034.32, "USD", 18973.41),
    @("XT0150", "63020051", "LKR", -29242804.73, "USD", -157042.07),
    @("UN0100", "63020602", "LKR", 13969885.84, "USD", 75022.22),
    @("UN0150", "63020621", "USD", -4375, "USD", -4375),
    @("XT0150", "63070001", "USD", -2532414, "USD", -2532414),
    @("XT0150", "63070002", "USD", 2853277.16, "USD", 2853277.16),
    @("UN0150", "63070012", "USD", 250851.46, "USD", 250851.46),
    @("UN0150", "63070501", "LKR", -301238357.4, "USD", -1617734.58),
    @("UN0150", "63070502", "LKR", 398347604.24, "USD", 2139238.51)
)

$r = 5
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Formula = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Column widths: A keeps its width (24); B/C/E are brand-new narrower
# columns; D and F keep the widths the old B/C columns used to have.
# ColumnWidth's character-width units run 5/6 higher than the raw stored
# <col width> value (Calibri 11 default-font padding), so back that off here
# to land on the exact target widths.
# ---------------------------------------------------------------------------
$pad = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 24 - $pad
$ws.Columns.Item(2).ColumnWidth = 12 - $pad
$ws.Columns.Item(3).ColumnWidth = 19 - $pad
$ws.Columns.Item(4).ColumnWidth = 22 - $pad
$ws.Columns.Item(5).ColumnWidth = 16 - $pad
$ws.Columns.Item(6).ColumnWidth = 26 - $pad

Write-Output "Summary sheet rebuilt to 6 columns"
